$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.967571021544131
$ws.Range("C2").Value = 0.1637334480314792
$ws.Range("D2").Value = 0.4266199014622885
$ws.Range("E2").Value = 0.1132740720607268
$ws.Range("G2").Value = 2.744180517121833
$ws.Range("H2").Value = 2.164845502587241
$ws.Range("I2").Value = 2.064196837229218
$ws.Range("J2").Value = 0.03841630471900359
$ws.Range("L2").Value = 0.5270728295534326
$ws.Range("M2").Value = 0.5186856003755764
$ws.Range("B3").Value = 1.890872322692871
$ws.Range("C3").Value = 0.1447541359414402
$ws.Range("D3").Value = 0.4268612994262497
$ws.Range("E3").Value = 0.1137379469430866
$ws.Range("G3").Value = 2.734431218647131
$ws.Range("H3").Value = 2.167547073505574
$ws.Range("I3").Value = 2.069047276777681
$ws.Range("J3").Value = 0.0369281932152461
$ws.Range("L3").Value = 0.5242597289030186
$ws.Range("M3").Value = 0.506099083042379
$ws.Range("B4").Value = 1.844769412694802
$ws.Range("C4").Value = 0.1331379115483742
$ws.Range("D4").Value = 0.4271686745170129
$ws.Range("E4").Value = 0.1140405009964129
$ws.Range("G4").Value = 2.729935943309243
$ws.Range("H4").Value = 2.170171481535164
$ws.Range("I4").Value = 2.072979670394673
$ws.Range("J4").Value = 0.03600117836106165
$ws.Range("L4").Value = 0.5227580552600273
$ws.Range("M4").Value = 0.4986281983259531
$ws.Range("B5").Value = 1.826231723026069
$ws.Range("C5").Value = 0.1284134402068275
$ws.Range("D5").Value = 0.4273339847342115
$ws.Range("E5").Value = 0.1141682619170927
$ws.Range("G5").Value = 2.728477983413086
$ws.Range("H5").Value = 2.171483475714922
$ws.Range("I5").Value = 2.074821666115312
$ws.Range("J5").Value = 0.03562004870457613
$ws.Range("L5").Value = 0.5222028807414603
$ws.Range("M5").Value = 0.4956485522005067
$ws.Range("B6").Value = 1.823168647764589
$ws.Range("C6").Value = 0.127629500047135
$ws.Range("D6").Value = 0.4273638543759901
$ws.Range("E6").Value = 0.1141897466104336
$ws.Range("G6").Value = 2.728258444428747
$ws.Range("H6").Value = 2.171715970964925
$ws.Range("I6").Value = 2.075141981617236
$ws.Range("J6").Value = 0.03555655876875008
$ws.Range("L6").Value = 0.522114124567068
$ws.Range("M6").Value = 0.4951577011303669
$ws.Range("B7").Value = 1.844518395449711
$ws.Range("C7").Value = 0.133074158422545
$ws.Range("D7").Value = 0.4271707417350896
$ws.Range("E7").Value = 0.114042205923917
$ws.Range("G7").Value = 2.729914768150962
$ws.Range("H7").Value = 2.170188193965515
$ws.Range("I7").Value = 2.073003542994272
$ws.Range("J7").Value = 0.03599605195458011
$ws.Range("L7").Value = 0.5227503380623375
$ws.Range("M7").Value = 0.4985877512647718
$ws.Range("B8").Value = 1.9409200619084
$ws.Range("C8").Value = 0.1571815991714232
$ws.Range("D8").Value = 0.4266701204891916
$ws.Range("E8").Value = 0.1134303408148538
$ws.Range("G8").Value = 2.740508935433382
$ws.Range("H8").Value = 2.165576408029267
$ws.Range("I8").Value = 2.065670961951696
$ws.Range("J8").Value = 0.03790595533883945
$ws.Range("L8").Value = 0.5260560775649026
$ws.Range("M8").Value = 0.5142924242991782
$ws.Range("B9").Value = 2.137809718738083
$ws.Range("C9").Value = 0.2047584342905111
$ws.Range("D9").Value = 0.4269502555145408
$ws.Range("E9").Value = 0.1123707872186719
$ws.Range("G9").Value = 2.773161873653407
$ws.Range("H9").Value = 2.164209672494621
$ws.Range("I9").Value = 2.058883895543389
$ws.Range("J9").Value = 0.04154659473463695
$ws.Range("L9").Value = 0.5343275230046913
$ws.Range("M9").Value = 0.5471286064829002
$ws.Range("B10").Value = 2.287251141630975
$ws.Range("C10").Value = 0.2399114540610583
$ws.Range("D10").Value = 0.4279246104307362
$ws.Range("E10").Value = 0.1116773236940176
$ws.Range("G10").Value = 2.804466934489596
$ws.Range("H10").Value = 2.16790948609804
$ws.Range("I10").Value = 2.058557035134243
$ws.Range("J10").Value = 0.0441590504693572
$ws.Range("L10").Value = 0.5414952820535319
$ws.Range("M10").Value = 0.572497379253214
$ws.Range("B11").Value = 2.356277654694907
$ws.Range("C11").Value = 0.2559498604182124
$ws.Range("D11").Value = 0.4285346564553976
$ws.Range("E11").Value = 0.1113801878686402
$ws.Range("G11").Value = 2.820313231063125
$ws.Range("H11").Value = 2.170619465028352
$ws.Range("I11").Value = 2.059427067895342
$ws.Range("J11").Value = 0.0453343276036513
$ws.Range("L11").Value = 0.5449930953750197
$ws.Range("M11").Value = 0.5843088976680306
$ws.Range("B12").Value = 2.382566286654708
$ws.Range("C12").Value = 0.2620301602799771
$ws.Range("D12").Value = 0.4287896356350132
$ws.Range("E12").Value = 0.1112702965670032
$ws.Range("G12").Value = 2.826545901132761
$ws.Range("H12").Value = 2.171793733407981
$ws.Range("I12").Value = 2.059903552596182
$ws.Range("J12").Value = 0.04577750718065943
$ws.Range("L12").Value = 0.5463517203715469
$ws.Range("M12").Value = 0.5888205765226999
$ws.Range("B13").Value = 2.376897904500879
$ws.Range("C13").Value = 0.2607203498005219
$ws.Range("D13").Value = 0.428733655610543
$ws.Range("E13").Value = 0.1112938468798315
$ws.Range("G13").Value = 2.825193246371242
$ws.Range("H13").Value = 2.171534242517566
$ws.Range("I13").Value = 2.059794385734207
$ws.Range("J13").Value = 0.04568214352290312
$ws.Range("L13").Value = 0.5460576008676981
$ws.Range("M13").Value = 0.5878471762373465
$ws.Range("B14").Value = 2.358437436801125
$ws.Range("C14").Value = 0.2564499512443774
$ws.Range("D14").Value = 0.4285551534951253
$ws.Range("E14").Value = 0.111371094426604
$ws.Range("G14").Value = 2.820821339864892
$ws.Range("H14").Value = 2.17071310324215
$ws.Range("I14").Value = 2.059463318978899
$ws.Range("J14").Value = 0.04537082567813044
$ws.Range("L14").Value = 0.5451041874107858
$ws.Range("M14").Value = 0.5846792967415126
$ws.Range("B15").Value = 2.347149361128686
$ws.Range("C15").Value = 0.2538351117224522
$ws.Range("D15").Value = 0.4284489365989259
$ws.Range("E15").Value = 0.1114187527499877
$ws.Range("G15").Value = 2.81817367143563
$ws.Range("H15").Value = 2.170229424905187
$ws.Range("I15").Value = 2.059279693490282
$ws.Range("J15").Value = 0.04517989119382193
$ws.Range("L15").Value = 0.5445246315344718
$ws.Range("M15").Value = 0.582743944557393
$ws.Range("B16").Value = 2.282761081973263
$ws.Range("C16").Value = 0.2388642654736941
$ws.Range("D16").Value = 0.4278880974539874
$ws.Range("E16").Value = 0.1116971103138331
$ws.Range("G16").Value = 2.803463746565342
$ws.Range("H16").Value = 2.167753079891952
$ws.Range("I16").Value = 2.058520721057931
$ws.Range("J16").Value = 0.04408198037862832
$ws.Range("L16").Value = 0.5412714618900907
$ws.Range("M16").Value = 0.5717309180995045
$ws.Range("B17").Value = 2.243528236013674
$ws.Range("C17").Value = 0.2296922942846322
$ws.Range("D17").Value = 0.4275867504956778
$ws.Range("E17").Value = 0.1118725616394567
$ws.Range("G17").Value = 2.794851703338395
$ws.Range("H17").Value = 2.166497215125332
$ws.Range("I17").Value = 2.05831641197625
$ws.Range("J17").Value = 0.04340509399493442
$ws.Range("L17").Value = 0.53933647258593
$ws.Range("M17").Value = 0.5650441719985437
$ws.Range("B18").Value = 2.2210609136892
$ws.Range("C18").Value = 0.2244212464074167
$ws.Range("D18").Value = 0.4274291254236999
$ws.Range("E18").Value = 0.111975201802573
$ws.Range("G18").Value = 2.79004932836051
$ws.Range("H18").Value = 2.165871509490302
$ws.Range("I18").Value = 2.058294753873213
$ws.Range("J18").Value = 0.04301452879622758
$ws.Range("L18").Value = 0.5382458422389362
$ws.Range("M18").Value = 0.5612236573689913
$ws.Range("B19").Value = 2.213470774507186
$ws.Range("C19").Value = 0.2226373167656845
$ws.Range("D19").Value = 0.4273784536540859
$ws.Range("E19").Value = 0.1120102505298184
$ws.Range("G19").Value = 2.78844923361757
$ws.Range("H19").Value = 2.165676241611322
$ws.Range("I19").Value = 2.058303867888696
$ws.Range("J19").Value = 0.04288207696669488
$ws.Range("L19").Value = 0.5378804085611506
$ws.Range("M19").Value = 0.5599344830762973
$ws.Range("B20").Value = 2.247694461644983
$ws.Range("C20").Value = 0.2306682071463797
$ws.Range("D20").Value = 0.427617204463985
$ws.Range("E20").Value = 0.1118537060485774
$ws.Range("G20").Value = 2.795752829195663
$ws.Range("H20").Value = 2.166620900003124
$ws.Range("I20").Value = 2.058328236592004
$ws.Range("J20").Value = 0.04347727771284937
$ws.Range("L20").Value = 0.5395401452076243
$ws.Range("M20").Value = 0.5657533463788127
$ws.Range("B21").Value = 2.363855661939283
$ws.Range("C21").Value = 0.2577040824187691
$ws.Range("D21").Value = 0.4286069335534393
$ws.Range("E21").Value = 0.1113483336964389
$ws.Range("G21").Value = 2.822099168517695
$ws.Range("H21").Value = 2.170950270367513
$ws.Range("I21").Value = 2.059556566916271
$ws.Range("J21").Value = 0.04546231795118416
$ws.Range("L21").Value = 0.5453833033464548
$ws.Range("M21").Value = 0.5856087242273205
$ws.Range("B22").Value = 2.440646615335766
$ws.Range("C22").Value = 0.2754138910188999
$ws.Range("D22").Value = 0.4293934745926009
$ws.Range("E22").Value = 0.1110333558404557
$ws.Range("G22").Value = 2.840670966602971
$ws.Range("H22").Value = 2.174642944945958
$ws.Range("I22").Value = 2.06121655909763
$ws.Range("J22").Value = 0.04674875117422062
$ws.Range("L22").Value = 0.5494007736374442
$ws.Range("M22").Value = 0.5988121456131665
$ws.Range("B23").Value = 2.399582113362044
$ws.Range("C23").Value = 0.2659581019322559
$ws.Range("D23").Value = 0.4289609065527316
$ws.Range("E23").Value = 0.1112000668163819
$ws.Range("G23").Value = 2.830634669357721
$ws.Range("H23").Value = 2.17259298395652
$ws.Range("I23").Value = 2.060251978901547
$ws.Range("J23").Value = 0.04606314956512847
$ws.Range("L23").Value = 0.547238407195664
$ws.Range("M23").Value = 0.5917445041205482
$ws.Range("B24").Value = 2.245810636384363
$ws.Range("C24").Value = 0.2302269905988794
$ws.Range("D24").Value = 0.4276033875597278
$ws.Range("E24").Value = 0.1118622251447423
$ws.Range("G24").Value = 2.795344966773911
$ws.Range("H24").Value = 2.166564682101665
$ws.Range("I24").Value = 2.058322592304421
$ws.Range("J24").Value = 0.04344464786057856
$ws.Range("L24").Value = 0.5394479968456807
$ws.Range("M24").Value = 0.5654326545495536
$ws.Range("B25").Value = 2.083705612056008
$ws.Range("C25").Value = 0.1918538473889839
$ws.Range("D25").Value = 0.4267394596334952
$ws.Range("E25").Value = 0.1126424599870797
$ws.Range("G25").Value = 2.76304917670123
$ws.Range("H25").Value = 2.163755089614966
$ws.Range("I25").Value = 2.059903973928321
$ws.Range("J25").Value = 0.04057275824267847
$ws.Range("L25").Value = 0.5318983550226619
$ws.Range("M25").Value = 0.5380271593023878
